$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'28.640.53"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.90%  '

$ws.Range('D3').Value = "'1.851.08"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.48%  '

$ws.Range('D5').Value = "'335.90"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.33%  '

$ws.Range('E6').Value = '  -0.91%  '

$ws.Range('D7').Value = "'0.4664"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.96%  '

$ws.Range('D8').Value = "'0.3917"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.36%  '

$ws.Range('D9').Value = "'0.07879"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.99%  '

$ws.Range('D10').Value = "'0.9836"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.41%  '

$ws.Range('D11').Value = "'22.24"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.16%  '

$ws.Range('D12').Value = "'1.839.35"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.16%  '

$ws.Range('D13').Value = "'5.850"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.28%  '

$ws.Range('D14').Value = "'7.026"
$ws.Range('D14').Style = 'Normal'

$ws.Range('D15').Value = "'0.06766"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.38%  '

$ws.Range('E16').Value = '  -0.96%  '

$ws.Range('D17').Value = "'87.69"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.05%  '

$ws.Range('E18').Value = '  -2.34%  '

$ws.Range('D19').Value = "'17.03"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.86%  '

$ws.Range('E20').Value = '  -0.94%  '

$ws.Range('D21').Value = "'28.631.22"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.96%  '

$ws.Range('D22').Value = "'5.417"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.35%  '

$ws.Range('D23').Value = "'11.29"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.89%  '

$ws.Range('D24').Value = "'2.123"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.19%  '

$ws.Range('D25').Value = "'2.081.51"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.34%  '

$ws.Range('D26').Value = "'153.48"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.72%  '

$ws.Range('D27').Value = "'6.287"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.77%  '

$ws.Range('D28').Value = "'19.43"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.88%  '

$ws.Range('D29').Value = "'2.019"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.53%  '

$ws.Range('D30').Value = "'117.64"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.28%  '

$ws.Range('D31').Value = "'0.9792"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.43%  '

$ws.Range('D32').Value = "'0.09455"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.62%  '

$ws.Range('D33').Value = "'5.378"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.09%  '

$ws.Range('D34').Value = "'3.505"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.49%  '

$ws.Range('E36').Value = '  -2.59%  '

$ws.Range('D37').Value = "'0.02196"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.68%  '

$ws.Range('D38').Value = "'1.160"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.60%  '

$ws.Range('D39').Value = "'0.5700"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.78%  '

$ws.Range('D40').Value = "'7.579"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.95%  '

$ws.Range('E41').Value = '  -5.59%  '

$ws.Range('D42').Value = "'0.1786"
$ws.Range('D42').Style = 'Normal'

$ws.Range('D43').Value = "'2.355"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.25%  '

$ws.Range('D44').Value = "'1.248"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.55%  '

$ws.Range('D45').Value = "'0.5381"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.14%  '

$ws.Range('D46').Value = "'11.80"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.34%  '

$ws.Range('D47').Value = "'0.07141"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.49%  '

$ws.Range('D48').Value = "'1.910"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.27%  '

$ws.Range('D49').Value = "'114.49"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.94%  '

$ws.Range('D50').Value = "'43.76"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.16%  '

$ws.Range('E51').Value = '  -0.96%  '
